# geração de análises seriais
# Reorders the tied category rows in the "max-arrecad" and "tx-sucesso"
# sheets (rows whose B value is a duplicate get a new A-column ordering).

$wb = $excel.ActiveWorkbook

# --- Sheet "max-arrecad": reorder categories within the two tied groups ---
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A2").Value = "questoes_genero"
$wsMax.Range("A3").Value = "fantasia"
$wsMax.Range("A4").Value = "fiq"
$wsMax.Range("A5").Value = "ficcao_cientifica"
$wsMax.Range("A6").Value = "humor"
$wsMax.Range("A7").Value = "folclore"
$wsMax.Range("A8").Value = "religiosidade"
$wsMax.Range("A9").Value = "terror"

$wsMax.Range("A13").Value = "angelo_agostini"
$wsMax.Range("A14").Value = "hqmix"

# --- Sheet "tx-sucesso": swap the two tied rows ---
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A8").Value = "erotismo"
$wsTx.Range("A9").Value = "questoes_genero"
